# Reorders the comma-separated "Recorded By" names/emails in column G
# into descending ordinal (case-sensitive) alphabetical order, e.g.
#   "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"

function OrdinalCompare($a, $b) {
    $la = $a.Length
    $lb = $b.Length
    $n = $la
    if ($lb -lt $n) { $n = $lb }
    for ($k = 0; $k -lt $n; $k++) {
        $ca = [int][char]$a[$k]
        $cb = [int][char]$b[$k]
        if ($ca -lt $cb) { return -1 }
        if ($ca -gt $cb) { return 1 }
    }
    if ($la -lt $lb) { return -1 }
    if ($la -gt $lb) { return 1 }
    return 0
}

function SortDescendingOrdinal($items) {
    $list = New-Object System.Collections.ArrayList
    foreach ($it in $items) { [void]$list.Add($it) }
    $cnt = $list.Count
    for ($m = 1; $m -lt $cnt; $m++) {
        $key = $list[$m]
        $j = $m - 1
        $continue = 1
        while ($j -ge 0 -and $continue -eq 1) {
            $cmp = OrdinalCompare $list[$j] $key
            if ($cmp -lt 0) {
                $list[$j+1] = $list[$j]
                $j = $j - 1
            } else {
                $continue = 0
            }
        }
        $list[$j+1] = $key
    }
    return $list
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ur = $ws.UsedRange
$lastRow = $ur.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cellAddr = "G" + $row
    $cell = $ws.Range($cellAddr)
    $val = $cell.Value2
    if ($val -ne $null -and $val -ne "") {
        $parts = $val.Split(",") | ForEach-Object { $_.Trim() }
        $sorted = SortDescendingOrdinal $parts
        $joined = [string]::Join(", ", $sorted)
        if ($joined -ne $val) {
            $cell.Value = $joined
        }
    }
}

Write-Host "Done reordering Recorded By column"
